# Apply the Student-Management workbook edit:
#  - sheet1 ("Trang_tính1"): update students in rows 3-5 (search/add/edit/delete
#    operations performed in the WinForms app), remove the old row 6 (Donald
#    Trump), move the selection to H6.
#  - add a new sheet "Trang_tính2" holding the summary stats (total students,
#    % male, % female) shown by the app's dashboard.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- sheet1: update existing student rows -------------------------------

# Row 3: Tap Can Binh (ID 5, Female, China)
$ws1.Range("A3").Value = "5"
$ws1.Range("B3").Value = "Tap"
$ws1.Range("C3").Value = "Can Binh"
$ws1.Range("D3").Value = 20459.905996087964
$ws1.Range("E3").Value = "Female"
$ws1.Range("F3").Value = "2399"
$ws1.Range("G3").Value = "China"
$ws1.Range("H3").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Tap can binh.jpeg"

# Row 4: Duong Yen Nhi (ID 3, Female, China)
$ws1.Range("A4").Value = "3"
$ws1.Range("B4").Value = "Duong"
$ws1.Range("C4").Value = "Yen Nhi"
$ws1.Range("D4").Value = 36957.63684953703
$ws1.Range("E4").Value = "Female"
$ws1.Range("F4").Value = "1579"
$ws1.Range("G4").Value = "China"
$ws1.Range("H4").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Tap can binh.jpeg"

# Row 5: Vladimir Putin (ID 2, Male, Russia)
$ws1.Range("A5").Value = "2"
$ws1.Range("B5").Value = "Vladimir"
$ws1.Range("C5").Value = "Putin"
$ws1.Range("D5").Value = 20827.638431944444
$ws1.Range("E5").Value = "Male"
$ws1.Range("F5").Value = "1933"
$ws1.Range("G5").Value = "Russia"
$ws1.Range("H5").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Putin.jpeg"

# Row 6 (Donald Trump) no longer exists after the edits in the app.
$ws1.Rows(6).Delete()

# Selection ends up on H6 in the saved workbook.
$ws1.Range("H6").Select()

# ---- add sheet2: summary dashboard ---------------------------------------

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Trang_tính2"

$ws2.Range("A1").Value = "Tổng sinh viên "
$ws2.Range("B1").Value = "Tỉ lệ nam %"
$ws2.Range("C1").Value = "Tỉ lệ nữ %"

$ws2.Range("A2").Value = 5
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 100

$ws2.Columns.Item(1).ColumnWidth = 16.608072916666668
$ws2.Columns.Item(2).ColumnWidth = 12.385416666666666
$ws2.Columns.Item(3).ColumnWidth = 12.276041666666666

$ws2.Range("E4:E5").Select()

# leave sheet1 as the active/selected tab
$ws1.Activate()
